$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for football_field preference, mirroring the style of H1
$ws.Range("I1").Value = "football_field"

# Copy the formatting from column H (H1:H26) onto column I (I1:I26) so that the
# new column matches the look of the existing preference columns.
$ws.Range("H1:H26").Copy()
$ws.Range("I1:I26").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Re-set the header text, since PasteSpecial(formats) does not touch values.
$ws.Range("I1").Value = "football_field"

# Fill in the preference values for football_field per space
$values = @(0.7, 0.7, 0, 0, 0.5, 0.5, 0, 0.8, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $values[$i]
}

$ws.Range("I26").Select() | Out-Null
